$d = $word.ActiveDocument

# Locate the paragraph that ends with "LOB1004: Cálculo II (Requisito fraco)".
# Immediately after it the document has three trailing paragraphs that must be
# removed: an empty paragraph, "Ver no Jupiter Salvar em pdf Salvar em docx",
# and the "© 2020 ..." footer line. A final empty paragraph and the
# page-break paragraph that follow those three must be kept untouched.
$anchorIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOB1004: Cálculo II (Requisito fraco)*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -ne $null) {
    $firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
    $lastToRemove = $d.Paragraphs.Item($anchorIndex + 3)

    # Sanity-check the paragraphs we are about to delete before touching them.
    if ($lastToRemove.Range.Text -like "*© 2020*") {
        $deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
        $deleteRange.Delete()
    }
}
